$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "الرصيد الحالي" (current balance) column H values refreshed for the new
# report timestamp: rows 7, 8, 10, 11 go from "1:0" -> "2:0"; row 9
# (QANDOVERAL) goes from "0:0" -> "1:0".
$ws.Range("H7").Value = "2:0"
$ws.Range("H8").Value = "2:0"
$ws.Range("H9").Value = "1:0"
$ws.Range("H10").Value = "2:0"
$ws.Range("H11").Value = "2:0"
